# Prueba.xlsx -> DataEntry rework ("termina proyecto dataentry desde excel para serenity")
#
# - Sheet "Prueba" becomes a small "Login" sheet (user/password sample row).
# - Sheet "Hoja1" is repurposed into a name/product/send table (2 stacked
#   copies of a 3-row block).
# - A brand-new "Hoja2" sheet is appended: the same table, 4 copies deep,
#   with the last two copies rendered in an underlined font.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Prueba" -> "Login"
# ---------------------------------------------------------------------
$login = $wb.Worksheets.Item(1)
$login.Name = "Login"
$login.Cells.Clear()

$login.Cells.Item(1,1).Value = "user"
$login.Cells.Item(1,2).Value = "password"
$login.Cells.Item(2,1).Value = "CamiloCh"
$login.Cells.Item(2,2).Value = "azxs1234"

[void]$login.Range("D9").Select()

# ---------------------------------------------------------------------
# 2) "Hoja1" gets rewritten with the name/product/send table (2 blocks)
# ---------------------------------------------------------------------
$hoja1 = $wb.Worksheets.Item(2)
$hoja1.Cells.Clear()

function Fill-Table($ws, $startRow, $blocks) {
    $r = $startRow
    $ws.Cells.Item($r,1).Value = "name"
    $ws.Cells.Item($r,2).Value = "product"
    $ws.Cells.Item($r,3).Value = "send"
    $r = $r + 1

    for ($b = 0; $b -lt $blocks; $b++) {
        $ws.Cells.Item($r,1).Value = "Camio"
        $ws.Cells.Item($r,2).Value = "TC"
        $ws.Cells.Item($r,3).Value = 5500
        $r = $r + 1

        $ws.Cells.Item($r,1).Value = "Diana"
        $ws.Cells.Item($r,2).Value = "TD"
        $ws.Cells.Item($r,3).Value = 698
        $r = $r + 1

        $ws.Cells.Item($r,1).Value = "Isaac"
        $ws.Cells.Item($r,2).Value = "credito"
        $ws.Cells.Item($r,3).Value = 310
        $r = $r + 1
    }
    return $r
}

[void](Fill-Table $hoja1 1 2)

[void]$hoja1.Range("A1:C7").Select()

# ---------------------------------------------------------------------
# 3) Brand-new "Hoja2" sheet appended after "Hoja1", same table x4 blocks,
#    with the 3rd and 4th blocks (rows 8-13) underlined.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$hoja2 = $wb.Worksheets.Add($null, $lastSheet)
$hoja2.Name = "Hoja2"

$nextRow = Fill-Table $hoja2 1 4

$underlineStart = $nextRow - 6
$underlineRange = $hoja2.Range($hoja2.Cells.Item($underlineStart,1), $hoja2.Cells.Item($nextRow - 1,3))
$underlineRange.Font.Underline = $true

[void]$hoja2.Range("B18").Select()

# Make "Hoja2" the active/visible tab, as in the saved workbook.
$hoja2.Activate()
